# "sorted the excel sheets"
#
# The results table on Sheet1 (A1:F16) gets re-ordered ascending by
# column E ("total avg Error: ") using Excel's Sort feature/ribbon
# ("Data > Sort"). The header row (row 1) is included in the sorted
# range - i.e. Header = xlNo - so it moves together with the data;
# since its column-E cell holds text ("total avg Error: ") it naturally
# sorts after all the numeric values (Excel's ascending sort order is
# numbers, then text), so it ends up last, on row 16.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$sortRange = $ws.Range("A1:F16")
$keyRange  = $ws.Range("E1:E16")

$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($keyRange, [Microsoft.Office.Interop.Excel.XlSortOn]::xlSortOnValues, [Microsoft.Office.Interop.Excel.XlSortOrder]::xlAscending, $null, [Microsoft.Office.Interop.Excel.XlSortDataOption]::xlSortNormal)

$ws.Sort.SetRange($sortRange)
$ws.Sort.Header = [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlNo
$ws.Sort.MatchCase = $false
$ws.Sort.Orientation = [Microsoft.Office.Interop.Excel.XlSortOrientation]::xlSortColumns
$ws.Sort.SortMethod = [Microsoft.Office.Interop.Excel.XlSortMethod]::xlPinYin
$ws.Sort.Apply()

# Columns A and E were also widened a bit so the (now longer/varied)
# sorted labels and error values stay fully visible.
$ws.Columns.Item(1).ColumnWidth = 43.8333333333333
$ws.Columns.Item(5).ColumnWidth = 23.3333333333333
